# Add a new "Decode" worksheet as a copy of the "Encode" worksheet, then
# tweak the numeric values (and a few cell styles that drifted along with
# them) to match the edited data, mirroring the commit "added decode and
# fixed encode".

$wb = $excel.ActiveWorkbook

$wsEncode = $wb.Worksheets.Item("Encode")

# Make sure Encode is the active sheet before duplicating it and select all
# of its cells (this matches what happens after someone copies the whole
# sheet to make a new one - the source sheet is left with a "select all"
# selection instead of its previous cursor position).
$wsEncode.Activate()
$wsEncode.Cells.Select()

# Duplicate "Encode" right after itself - this carries over all formatting,
# column widths, shared styles and the worksheet's own data.
$wsEncode.Copy([System.Reflection.Missing]::Value, $wsEncode)

$wsDecode = $wb.Worksheets.Item($wsEncode.Index + 1)
$wsDecode.Name = "Decode"

# Re-create the Excel Table ("ListObject") that lived on Encode - the
# worksheet copy itself does not bring along the table definition.
$lo = $wsDecode.ListObjects.Add(1, $wsDecode.Range("A1:I20"), [System.Reflection.Missing]::Value, 1)
$lo.Name = "Tabelle145"
$lo.TableStyle = "TableStyleMedium4"

# ---- Block 1 (rows 2-7): Konstant -------------------------------------
$wsDecode.Range("D2").Value = 6000
$wsDecode.Range("E2").Value = 48
$wsDecode.Range("H2").Value = 8

# Only columns D, E and H pick up row 2's cell style further down the
# block (F and G keep their own formatting), so copy each column on its
# own rather than the whole D:H block.
$wsDecode.Range("D2").Copy()
$wsDecode.Range("D3:D7").PasteSpecial(-4122)
$wsDecode.Range("E2").Copy()
$wsDecode.Range("E3:E7").PasteSpecial(-4122)
$wsDecode.Range("H2").Copy()
$wsDecode.Range("H3:H7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsDecode.Range("D3:D7").Value = 6000
$wsDecode.Range("E3:E7").Value = 48
$wsDecode.Range("H3:H7").Value = 8

# ---- Block 2 (rows 8-13): One-Fifth -----------------------------------
$wsDecode.Range("H2").Copy()
$wsDecode.Range("H8:H13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsDecode.Range("D8:D13").Value = 2500
$wsDecode.Range("E8:E13").Value = 50
$wsDecode.Range("H8:H13").Value = 8

# ---- Block 3 (rows 14-19): linear fallend -----------------------------
$wsDecode.Range("D14:D19").Value = 3500
$wsDecode.Range("E14:E19").Value = 48

# ---- Row 20: Keine Rekombination summary row --------------------------
$wsDecode.Range("D20").Value = 5500
$wsDecode.Range("H20").Value = 6

# Leave the new sheet's selection somewhere in the data, matching the
# author's last edit position, and make it the active tab.
$wsDecode.Range("D11").Select()
$wsDecode.Activate()
